# Update "想去人数" (interested-count) figures in column F across the four
# sheets of the workbook, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 249
$ws.Range("F6").Value = 891
$ws.Range("F7").Value = 460
$ws.Range("F9").Value = 2168
$ws.Range("F10").Value = 622
$ws.Range("F11").Value = 285
$ws.Range("F13").Value = 1062
$ws.Range("F15").Value = 2195
$ws.Range("F16").Value = 652
$ws.Range("F17").Value = 12558
$ws.Range("F18").Value = 1233
$ws.Range("F19").Value = 8
$ws.Range("F20").Value = 555
$ws.Range("F21").Value = 126
$ws.Range("F22").Value = 21
$ws.Range("F23").Value = 138
$ws.Range("F24").Value = 39
$ws.Range("F25").Value = 262
$ws.Range("F27").Value = 3

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 18
$ws.Range("F9").Value = 148
$ws.Range("F11").Value = 82
$ws.Range("F12").Value = 58

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5695
$ws.Range("F4").Value = 465

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5695
$ws.Range("F5").Value = 465
$ws.Range("F9").Value = 891
$ws.Range("F11").Value = 460
$ws.Range("F12").Value = 18
$ws.Range("F13").Value = 2168
$ws.Range("F14").Value = 622
$ws.Range("F15").Value = 285
$ws.Range("F19").Value = 1062
$ws.Range("F22").Value = 148
$ws.Range("F24").Value = 2195
$ws.Range("F25").Value = 652
$ws.Range("F26").Value = 82
$ws.Range("F27").Value = 58
$ws.Range("F28").Value = 1233
$ws.Range("F29").Value = 8
$ws.Range("F30").Value = 555
$ws.Range("F31").Value = 126
$ws.Range("F32").Value = 21
$ws.Range("F33").Value = 138
$ws.Range("F35").Value = 39
$ws.Range("F38").Value = 262
$ws.Range("F44").Value = 3

$wb.Save()
